$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "b.md" file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 15:19:16"

# --- zh-cn sheet: row 3 is the "b.md" file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "False" is typed literally as text (not boolean) in the source data; a
# plain .Value assignment of the word False/True gets auto-coerced to a
# real COM boolean, so copy the existing text "False" cell (O3) over
# instead - Range.Copy preserves the underlying text cell type/style.
$wsZhCn.Range("O3").Copy($wsZhCn.Range("F3"))
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-06 15:18:59"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a72fcfa1e9d351cfe847c3a15ce311c9913c1823/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad10ce155b1feb41622702d73387599dd2d56ada/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is the "b.md" file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("O3").Copy($wsDeDe.Range("F3"))
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-06 15:19:16"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a72fcfa1e9d351cfe847c3a15ce311c9913c1823/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad10ce155b1feb41622702d73387599dd2d56ada/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
